# "after Haunting Mars 5"
# Update the character sheet from the "Synth (RR C Flexi-Skin)" morph to the
# "Splicer" morph, tweaking stats, traits and adding the new "Sex Appeal" skill.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("F")

# --- Morph name / description -------------------------------------------
$ws.Range("A4").Value2 = "Splicer"
$ws.Range("A7").Value2 = "Mesh Inserts, Cortical Stack, Basic Biomods, Handsome +2, Intersex, Pheromones"

# --- New "LINE 2:" helper cells (character description) ------------------
$ws.Range("I2").Value2 = "LINE 2:"
$ws.Range("J2").Value2 = "blond, 176 cm, 63 kg"
$ws.Range("K2").ClearContents()

# --- Base stat values (row 4) and morph stat bonuses (row 5) -------------
$ws.Range("C4").Value2 = 11
$ws.Range("D4").Value2 = 11
$ws.Range("E4").Value2 = 11
$ws.Range("F4").Value2 = 11
$ws.Range("G4").Value2 = 11
$ws.Range("H4").Value2 = 12

$ws.Range("E5").Value2 = 1
$ws.Range("F5").ClearContents()
$ws.Range("G5").Value2 = 1

# --- New skill row 32: Sex Appeal (based on HT) ---------------------------
$ws.Range("A32").Value2 = "Sex Appeal"
$ws.Range("B32").Value2 = "HT"
$ws.Range("C32").Value2 = -1
$ws.Range("D32").Value2 = 2
$ws.Range("E32").Formula = '=H32+C32+D32'
$ws.Range("H32").Formula = '=INDEX($C$3:$H$3,(MATCH($B32,$2:$2,0)-2))'

# --- Output sheet: LINE 1 now also appends the new LINE 2 description ----
$wsOut = $wb.Worksheets.Item("OUTPUT")
$wsOut.Range("A3").Formula = '=UPPER(F!A4)&": "&(F!J2)'

# --- Update selections to match the saved file ---------------------------
$ws.Range("D34").Select()
$wsOut2 = $wb.Worksheets.Item("OUTPUT")
$wsOut2.Range("A4").Select()

$wb.Application.Calculate()
